$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $val) {
    $ws.Range($addr).Value = "" + $val
    $ws.Range($addr).Style = "Normal"
    $ws.Range($addr).Font.Name = "Andale WT"
    $ws.Range($addr).Font.Size = 10
    $ws.Range($addr).HorizontalAlignment = -4152
    $ws.Range($addr).VerticalAlignment = -4108
}

# Header text updates (Volume Number, Report week dates)
$ws.Range("A8").Value = "Volume 30   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/16/2023  Through  1/22/2023"

# Simple numeric value updates (style unchanged)
$ws.Range("H14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("G15").Value = 2
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = 150
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 9
$ws.Range("J16").Value = 8
$ws.Range("K16").Value = 12.5
$ws.Range("M16").Value = 12.5
$ws.Range("N16").Value = -75.675675675675
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -50
$ws.Range("F17").Value = 17
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -32
$ws.Range("I17").Value = 12
$ws.Range("J17").Value = 22
$ws.Range("K17").Value = -45.454545454545
$ws.Range("L17").Value = 100
$ws.Range("M17").Value = 9.090909090909
$ws.Range("N17").Value = 20
$ws.Range("D18").Value = 1
$ws.Range("G18").Value = 6
$ws.Range("J18").Value = 5
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 60
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = 52.941176470588
$ws.Range("I19").Value = 19
$ws.Range("J19").Value = 13
$ws.Range("K19").Value = 46.153846153846
$ws.Range("L19").Value = 137.5
$ws.Range("M19").Value = 72.727272727272
$ws.Range("N19").Value = -13.636363636363
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -25
$ws.Range("F20").Value = 8
$ws.Range("G20").Value = 10
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 9
$ws.Range("K20").Value = -33.333333333333
$ws.Range("M20").Value = -53.846153846153
$ws.Range("N20").Value = -95.683453237410
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 18.75
$ws.Range("F21").Value = 62
$ws.Range("G21").Value = 69
$ws.Range("H21").Value = -10.144927536231
$ws.Range("I21").Value = 47
$ws.Range("J21").Value = 59
$ws.Range("K21").Value = -20.338983050847
$ws.Range("L21").Value = 74.074074074074
$ws.Range("M21").Value = -22.950819672131
$ws.Range("N21").Value = -81.274900398406
$ws.Range("C23").Value = 2
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 9
$ws.Range("G23").Value = 6
$ws.Range("H23").Value = 50
$ws.Range("I23").Value = 6
$ws.Range("J23").Value = 3
$ws.Range("L23").Value = 200
$ws.Range("M23").Value = 500
$ws.Range("C24").Value = 7
$ws.Range("E24").Value = -41.666666666666
$ws.Range("F24").Value = 55
$ws.Range("G24").Value = 40
$ws.Range("H24").Value = 37.5
$ws.Range("I24").Value = 35
$ws.Range("J24").Value = 29
$ws.Range("K24").Value = 20.689655172413
$ws.Range("L24").Value = -2.777777777777
$ws.Range("M24").Value = 16.666666666666
$ws.Range("C25").Value = 4
$ws.Range("E25").Value = -55.555555555555
$ws.Range("F25").Value = 20
$ws.Range("G25").Value = 31
$ws.Range("H25").Value = -35.483870967741
$ws.Range("I25").Value = 15
$ws.Range("J25").Value = 23
$ws.Range("K25").Value = -34.782608695652
$ws.Range("L25").Value = 66.666666666666
$ws.Range("M25").Value = -44.444444444444
$ws.Range("G26").Value = 2
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -50
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = -50
$ws.Range("L27").Value = 100
$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -80
$ws.Range("J28").Value = 5
$ws.Range("K28").Value = -80
$ws.Range("L28").Value = -50
$ws.Range("M28").Value = -50
$ws.Range("N28").Value = -66.666666666666
$ws.Range("E29").Value = 0
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = -50
$ws.Range("J29").Value = 2
$ws.Range("K29").Value = -50
$ws.Range("L29").Value = -50
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -66.666666666666

# Numeric updates requiring style change (text -> number)
$ws.Range("C14").Value = 1
$ws.Range("C14").NumberFormat = "#,##0"
$ws.Range("F14").Value = 1
$ws.Range("F14").NumberFormat = "#,##0"
$ws.Range("I14").Value = 1
$ws.Range("I14").NumberFormat = "#,##0"
$ws.Range("L16").Value = 800
$ws.Range("L16").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("C27").Value = 1
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("F27").Value = 2
$ws.Range("F27").NumberFormat = "#,##0"
$ws.Range("I27").Value = 2
$ws.Range("I27").NumberFormat = "#,##0"
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("F28").Value = 1
$ws.Range("F28").NumberFormat = "#,##0"
$ws.Range("I28").Value = 1
$ws.Range("I28").NumberFormat = "#,##0"
$ws.Range("C29").Value = 1
$ws.Range("C29").NumberFormat = "#,##0"
$ws.Range("F29").Value = 1
$ws.Range("F29").NumberFormat = "#,##0"
$ws.Range("I29").Value = 1
$ws.Range("I29").NumberFormat = "#,##0"

# Text updates (number -> text)
Set-TextCell "D14" "0"
Set-TextCell "E14" "***.*"
Set-TextCell "D15" "0"
Set-TextCell "E15" "***.*"
Set-TextCell "D26" "0"
Set-TextCell "E26" "***.*"
